# Adicao painel de administracao
# Updates the Funcionarios sheet: removes the shared "Booker@1010" password
# hyperlink from column B, re-points several employees' department/manager
# (columns E/G) at real people instead of the placeholder "sistema", and
# appends three new admin ("Socio") accounts as rows 69-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Funcionarios")

# ---------------------------------------------------------------------
# 1) Clear the "senha" (password) column B for every data row (2-68).
#    The cells keep their existing (Hiperlink) style but lose both their
#    text value and their hyperlink.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 68; $r++) {
    $ws.Cells.Item($r, 2).ClearContents()
}

# Hyperlinks.Delete() on this host clears the whole worksheet collection,
# so wipe it once and re-add only the link(s) that must survive.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C27"), "mailto:joaovitor.alves@bookerbrasil.com") | Out-Null

# ---------------------------------------------------------------------
# 2) Point specific employees' departamento (E) / nome_gestor (G) at the
#    real department + manager instead of the "sistema" placeholder.
# ---------------------------------------------------------------------
$deptGestor = @{
    13 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    15 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    27 = @("TI", "RUDGE RODRIGUES")
    30 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    31 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    48 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    53 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    57 = @("DEPARTAMENTO PESSOAL", "KARINE CASTRO")
    62 = @("DEPARTAMENTO PESSOAL", "REINALDO RODRIGUES DAMASCENA")
    65 = @("TI", "RUDGE RODRIGUES")
    68 = @("TI", "RUDGE RODRIGUES")
}

foreach ($r in $deptGestor.Keys) {
    $vals = $deptGestor[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 7).Value = $vals[1]
}

# ---------------------------------------------------------------------
# 3) Append three new "Socio" (admin) accounts as rows 69-71.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 69; A = "karine.castro";   C = "karine.castro@bookerbrasil.com";   D = "KARINE CASTRO";   E = "DEPARTAMENTO PESSOAL"; F = "Sócio" }
    @{ Row = 70; A = "rudge.rodrigues"; C = "rudge.rodrigues@bookerbrasil.com"; D = "RUDGE RODRIGUES"; E = "CEO";                  F = "Sócio" }
    @{ Row = 71; A = "leandro.santana"; C = "leandro.santana@bookerbrasil.com"; D = "LEANDRO SANTANA"; E = "CONTABIL";             F = "Sócio" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 3).Value = $nr.C
    $ws.Cells.Item($r, 4).Value = $nr.D
    $ws.Cells.Item($r, 5).Value = $nr.E
    $ws.Cells.Item($r, 6).Value = $nr.F
    $ws.Cells.Item($r, 8).Value = "Sim"

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 3), "mailto:" + $nr.C) | Out-Null
    $ws.Cells.Item($r, 3).Style = "Hiperlink"
}

# ---------------------------------------------------------------------
# 4) Restore the header cell / selection the author left active.
# ---------------------------------------------------------------------
$ws.Range("F65").Select()
